$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Copy()
$ws.Range("A10:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value2 = 43382
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Meeting"

$ws.Range("A11").Value2 = 43386
$ws.Range("B11").Value = 1.5
$ws.Range("C11").Value = "Coderen fietssimulatie drag + timeslots"

$ws.Range("A12").Value2 = 43387
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Keras/TF werkende krijgen"

$ws.Range("A13").Value2 = 43387
$ws.Range("B13").Value = 1.5
$ws.Range("C13").Value = "Keras bekijken/leren"

$ws.Range("B14").Select()
